# Auto-generated Excel COM-interop edit script
# Applies numeric "want-to-go count" corrections across sheets and
# re-syncs the "全部类型" (all-types) aggregate sheet rows 27-37
# with the latest per-category event listings.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, [string]$addr, [string]$text)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

function Set-NumberCell {
    param($ws, [string]$addr, $number)
    $ws.Range($addr).Value = $number
}

# --- Sheet "展览": bump 想去人数 (want-to-go counts) in column F ---
$wsExpo = $wb.Worksheets.Item("展览")
Set-NumberCell $wsExpo "F6" 81
Set-NumberCell $wsExpo "F7" 4522
Set-NumberCell $wsExpo "F10" 2642
Set-NumberCell $wsExpo "F14" 1686
Set-NumberCell $wsExpo "F15" 699
Set-NumberCell $wsExpo "F16" 447
Set-NumberCell $wsExpo "F17" 166
Set-NumberCell $wsExpo "F23" 497
Set-NumberCell $wsExpo "F24" 36
Set-NumberCell $wsExpo "F32" 1274
Set-NumberCell $wsExpo "F33" 242
Set-NumberCell $wsExpo "F35" 1317
Set-NumberCell $wsExpo "F37" 333
Set-NumberCell $wsExpo "F43" 719
Set-NumberCell $wsExpo "F44" 1400
Set-NumberCell $wsExpo "F45" 156

# --- Sheet "全部类型": bump 想去人数 (want-to-go counts) in column F ---
$wsAll = $wb.Worksheets.Item("全部类型")
Set-NumberCell $wsAll "F4" 81
Set-NumberCell $wsAll "F5" 4522
Set-NumberCell $wsAll "F7" 2642
Set-NumberCell $wsAll "F8" 1686
Set-NumberCell $wsAll "F11" 699
Set-NumberCell $wsAll "F12" 447
Set-NumberCell $wsAll "F13" 166
Set-NumberCell $wsAll "F19" 497
Set-NumberCell $wsAll "F20" 36
Set-NumberCell $wsAll "F43" 719
Set-NumberCell $wsAll "F44" 1400
Set-NumberCell $wsAll "F46" 156

# --- Sheet "全部类型": resync rows 27-37 with latest event listing ---
# (new event inserted, two stale/removed events dropped, remaining rows
#  shift to reflect current chronological listing)
# row 27
Set-TextCell $wsAll "B27" "2024-04-27"
Set-TextCell $wsAll "C27" "杭州·原神x崩铁周年特典only"
Set-TextCell $wsAll "D27" "康候圣街99号 顺丰创新中心"
Set-TextCell $wsAll "E27" "2024.04.27 09:00-04.27 17:00"
Set-NumberCell $wsAll "F27" 15
Set-NumberCell $wsAll "G27" 58
Set-TextCell $wsAll "H27" "https://show.bilibili.com/platform/detail.html?id=83495"
Set-TextCell $wsAll "I27" "//i1.hdslb.com/bfs/openplatform/202403/LqG2INbt1711358703138.png"

# row 28
Set-TextCell $wsAll "B28" "2024-05-01"
Set-TextCell $wsAll "C28" "杭州·与梦回望动漫游戏展"
Set-TextCell $wsAll "D28" "沈半路171号 T-Car杭州汽车文化主题公园"
Set-TextCell $wsAll "E28" "2024.05.01 10:00-05.02 17:00"
Set-NumberCell $wsAll "F28" 461
Set-NumberCell $wsAll "G28" 70
Set-TextCell $wsAll "H28" "https://show.bilibili.com/platform/detail.html?id=82725"
Set-TextCell $wsAll "I28" "//i0.hdslb.com/bfs/openplatform/202403/lt13shal1710228931298.jpeg"

# row 29
Set-TextCell $wsAll "B29" "2024-05-01"
Set-TextCell $wsAll "C29" "杭州·第37届 中二病 原神x星穹only"
Set-TextCell $wsAll "D29" "康候圣街99号 顺丰创新中心"
Set-TextCell $wsAll "E29" "2024.05.01 10:30-05.02 17:00"
Set-NumberCell $wsAll "F29" 1640
Set-NumberCell $wsAll "G29" 60
Set-TextCell $wsAll "H29" "https://show.bilibili.com/platform/detail.html?id=82700"
Set-TextCell $wsAll "I29" "//i1.hdslb.com/bfs/openplatform/202403/Kb75MESZ1710215541381.jpeg"

# row 30
Set-TextCell $wsAll "B30" "2024-05-01"
Set-TextCell $wsAll "C30" "杭州·第7届YH樱花动漫游戏文化节"
Set-TextCell $wsAll "D30" "德胜东路2539号 梦马汽车小镇"
Set-TextCell $wsAll "E30" "2024.05.01 10:00-05.02 17:00"
Set-NumberCell $wsAll "F30" 1274
Set-NumberCell $wsAll "G30" 70
Set-TextCell $wsAll "H30" "https://show.bilibili.com/platform/detail.html?id=82828"
Set-TextCell $wsAll "I30" "//i1.hdslb.com/bfs/openplatform/202403/Kd0niodt1710905544733.jpeg"

# row 31
Set-TextCell $wsAll "B31" "2024-05-01"
Set-TextCell $wsAll "C31" "杭州·第7届YH樱花漫展-SVIP嘉宾前排票"
Set-TextCell $wsAll "D31" "德胜东路2539号 梦马汽车小镇"
Set-TextCell $wsAll "E31" "2024.05.01 10:00-05.02 17:00"
Set-NumberCell $wsAll "F31" 242
Set-NumberCell $wsAll "G31" 168
Set-TextCell $wsAll "H31" "https://show.bilibili.com/platform/detail.html?id=83267"
Set-TextCell $wsAll "I31" "//i1.hdslb.com/bfs/openplatform/202403/DgmIZ6G71711357279757.jpeg"

# row 32
Set-TextCell $wsAll "B32" "2024-05-01"
Set-TextCell $wsAll "C32" "杭州·第7届YH樱花漫展-配音演员紫枫儿内场票"
Set-TextCell $wsAll "D32" "德胜东路2539号 梦马汽车小镇"
Set-TextCell $wsAll "E32" "2024.05.01 10:00-05.01 17:00"
Set-NumberCell $wsAll "F32" 33
Set-NumberCell $wsAll "G32" 98
Set-TextCell $wsAll "H32" "https://show.bilibili.com/platform/detail.html?id=83331"
Set-TextCell $wsAll "I32" "//i0.hdslb.com/bfs/openplatform/202403/h5ilz3SA1711351453471.jpeg"

# row 33
Set-TextCell $wsAll "B33" "2024-05-02"
Set-TextCell $wsAll "C33" "杭州·第四届华盟动漫次元嘉年华"
Set-TextCell $wsAll "D33" "创意路1号 中国智谷富春园区"
Set-TextCell $wsAll "E33" "2024.05.02 10:00-05.03 17:00"
Set-NumberCell $wsAll "F33" 2176
Set-NumberCell $wsAll "G33" 58
Set-TextCell $wsAll "H33" "https://show.bilibili.com/platform/detail.html?id=82465"
Set-TextCell $wsAll "I33" "//i0.hdslb.com/bfs/openplatform/202403/4XHyqi3D1709780326858.jpeg"

# row 34
Set-TextCell $wsAll "B34" "2024-05-02"
Set-TextCell $wsAll "C34" "杭州·造梦探险家Porject6野蛮冲撞——第五人格ONLY"
Set-TextCell $wsAll "D34" "欢西路1号 天都城酒店"
Set-TextCell $wsAll "E34" "2024.05.02 10:00-05.02 22:00"
Set-NumberCell $wsAll "F34" 333
Set-NumberCell $wsAll "G34" 28
Set-TextCell $wsAll "H34" "https://show.bilibili.com/platform/detail.html?id=82851"
Set-TextCell $wsAll "I34" "//i1.hdslb.com/bfs/openplatform/202403/a7IYN66u1711441126355.png"

# row 35
Set-TextCell $wsAll "B35" "2024-05-12"
Set-TextCell $wsAll "C35" "杭州·《卡农》永恒经典名曲音乐会"
Set-TextCell $wsAll "D35" "武林路77号 浙江省文化馆小剧场（原群艺馆小剧场）"
Set-TextCell $wsAll "E35" "2024.05.12 14:00-05.12 15:30"
Set-NumberCell $wsAll "F35" 2
Set-NumberCell $wsAll "G35" 100
Set-TextCell $wsAll "H35" "https://show.bilibili.com/platform/detail.html?id=83176"
Set-TextCell $wsAll "I35" "//i0.hdslb.com/bfs/openplatform/202403/gLrSkh0O1711013683966.jpeg"

# row 36
Set-TextCell $wsAll "B36" "2024-05-12"
Set-TextCell $wsAll "C36" "杭州·奇迹の闪耀 「UP!」巡回动漫演唱会"
Set-TextCell $wsAll "D36" "东坡路10号 杭州东坡大剧院"
Set-TextCell $wsAll "E36" "2024.05.12 19:30-05.12 21:30"
Set-NumberCell $wsAll "F36" 12
Set-NumberCell $wsAll "G36" 126
Set-TextCell $wsAll "H36" "https://show.bilibili.com/platform/detail.html?id=82452"
Set-TextCell $wsAll "I36" "//i1.hdslb.com/bfs/openplatform/202403/HvxHPz981709707512970.jpeg"

# row 37
Set-TextCell $wsAll "B37" "2024-05-18"
Set-TextCell $wsAll "C37" "杭州·Jo迪"
Set-TextCell $wsAll "D37" "萧杭路28号 格拉斯club"
Set-TextCell $wsAll "E37" "2024.05.18 13:00-05.18 19:00"
Set-NumberCell $wsAll "F37" 18
Set-NumberCell $wsAll "G37" 198
Set-TextCell $wsAll "H37" "https://show.bilibili.com/platform/detail.html?id=83008"
Set-TextCell $wsAll "I37" "//i1.hdslb.com/bfs/openplatform/202403/AEtl5BHN1711015003341.jpeg"

